$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = ($cols[$i] + "_FV2310")
}
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = ($cols[$i] + "_FV2404")
}
